$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers updated with new timestamps)
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555566310518"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555599904544"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555599914553"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555600554547"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555601344533"

# Sheet 1: GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555565990508.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555566130502.csv"
$ws1.Range("B4").Value = "go_stims-16512555566150572.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555566290507.csv"

# Sheet 2: NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16512555587878757.csv"
$ws2.Range("B3").Value = "OB-16512555582138758.csv"
$ws2.Range("B4").Value = "TB-16512555599714537.csv"
$ws2.Range("B5").Value = "ZB-match_3-16512555571451917.csv"
$ws2.Range("B6").Value = "ZB-match_9-16512555571951911.csv"
$ws2.Range("B7").Value = "TB-16512555586228735.csv"
$ws2.Range("B8").Value = "OB-16512555582668757.csv"
$ws2.Range("B9").Value = "OB-1651255558575874.csv"
$ws2.Range("B10").Value = "ZB-match_0-16512555572912407.csv"

# Sheet 3: RS - no cell content changes, only sheet name updated above

# Sheet 4: TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555600224576.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555599994588.csv"
$ws4.Range("B4").Value = "MM_stims-16512555600384576.csv"
$ws4.Range("B5").Value = "ZM_stims-165125556002346.csv"
$ws4.Range("B6").Value = "MM_stims-16512555600544565.csv"
$ws4.Range("B7").Value = "ZM_stims-1651255560039458.csv"

# Sheet 5: vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512555600874565.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555601184535.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555601024542.csv"
$ws5.Range("B5").Value = "SAT_stims-1651255560061455.csv"
